# Generate Report for Handback
#
# The localization-status report tracks two source files:
#   63a398a4-7e50-4038-a6ee-4d0ede53b8bc.md  (still "Ready for handoff")
#   d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md  (now handed back)
#
# This run records the handback for d12631b3-...md: its status flips to
# "Handed back: in sync with en-US", its handback file/datetime columns
# get populated, and (because the report re-sorts by latest activity) its
# row moves above the still-pending 63a398a4-...md row on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

# Row 2 now reflects d12631b3-...md (handed back)
$wsOverview.Range("B2").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value2 = "2016-03-22 16:40:47"

# Row 3 now reflects 63a398a4-...md (still ready for handoff)
$wsOverview.Range("B3").Value2 = "Ready for handoff"
$wsOverview.Range("C3").Value2 = "Ready for handoff"
$wsOverview.Range("D3").Value2 = "2016-03-22 16:40:12"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/66e1ccdebd1552c8c00808215c9ec71d1eed7b89/e2e/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/08c1d6435e5a6bfbe274169f4ece73b3e1834d4f/e2e/63a398a4-7e50-4038-a6ee-4d0ede53b8bc.md", "", "", "63a398a4-7e50-4038-a6ee-4d0ede53b8bc.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": per-language handoff/handback detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Delete()

# Row 2 now reflects d12631b3-...md (handed back)
$wsZhCn.Range("B2").Value2 = ".md"
$wsZhCn.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsZhCn.Range("E2").Value2 = "2016-03-22 16:40:43"
$wsZhCn.Range("H2").Value2 = "2016-03-22 16:41:33"
$wsZhCn.Range("J2").Value2 = "Include"

# Row 3 now reflects 63a398a4-...md (still ready for handoff)
$wsZhCn.Range("B3").Value2 = ".md"
$wsZhCn.Range("C3").Value2 = "Ready for handoff"
$wsZhCn.Range("E3").Value2 = "2016-03-22 16:40:07"
$wsZhCn.Range("H3").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("J3").Value2 = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/66e1ccdebd1552c8c00808215c9ec71d1eed7b89/e2e/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/743e257967a1356ccf74872997faaf23d632702f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.1b507a749cb9fa63b0ac0e39e0c1fc790417cdb0.zh-cn.xlf", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.1b507a749cb9fa63b0ac0e39e0c1fc790417cdb0.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/66e1ccdebd1552c8c00808215c9ec71d1eed7b89/e2e/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/743e257967a1356ccf74872997faaf23d632702f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.1b507a749cb9fa63b0ac0e39e0c1fc790417cdb0.zh-cn.xlf", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.1b507a749cb9fa63b0ac0e39e0c1fc790417cdb0.zh-cn.xlf") | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/08c1d6435e5a6bfbe274169f4ece73b3e1834d4f/e2e/63a398a4-7e50-4038-a6ee-4d0ede53b8bc.md", "", "", "63a398a4-7e50-4038-a6ee-4d0ede53b8bc.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/55a036f50a2f7db0686c3fe5712552726535c9b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/63a398a4-7e50-4038-a6ee-4d0ede53b8bc.3026d92881defd4f39bed6999f1fc2987b39984b.zh-cn.xlf", "", "", "63a398a4-7e50-4038-a6ee-4d0ede53b8bc.3026d92881defd4f39bed6999f1fc2987b39984b.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": per-language handoff/handback detail
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Delete()

# Row 2 now reflects d12631b3-...md (handed back)
$wsDeDe.Range("B2").Value2 = ".md"
$wsDeDe.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsDeDe.Range("E2").Value2 = "2016-03-22 16:40:47"
$wsDeDe.Range("H2").Value2 = "2016-03-22 16:41:39"
$wsDeDe.Range("J2").Value2 = "Include"

# Row 3 now reflects 63a398a4-...md (still ready for handoff)
$wsDeDe.Range("B3").Value2 = ".md"
$wsDeDe.Range("C3").Value2 = "Ready for handoff"
$wsDeDe.Range("E3").Value2 = "2016-03-22 16:40:12"
$wsDeDe.Range("H3").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("J3").Value2 = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/66e1ccdebd1552c8c00808215c9ec71d1eed7b89/e2e/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbdca95033c58f6b18be8a46fd5a2ff859c7e90d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.1b507a749cb9fa63b0ac0e39e0c1fc790417cdb0.de-de.xlf", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.1b507a749cb9fa63b0ac0e39e0c1fc790417cdb0.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/66e1ccdebd1552c8c00808215c9ec71d1eed7b89/e2e/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbdca95033c58f6b18be8a46fd5a2ff859c7e90d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.1b507a749cb9fa63b0ac0e39e0c1fc790417cdb0.de-de.xlf", "", "", "d12631b3-8dcb-4c3d-9ce6-fe3dbd404716.1b507a749cb9fa63b0ac0e39e0c1fc790417cdb0.de-de.xlf") | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/08c1d6435e5a6bfbe274169f4ece73b3e1834d4f/e2e/63a398a4-7e50-4038-a6ee-4d0ede53b8bc.md", "", "", "63a398a4-7e50-4038-a6ee-4d0ede53b8bc.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16c75f4209c73a4375aad088b27195c2506b86bc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/63a398a4-7e50-4038-a6ee-4d0ede53b8bc.3026d92881defd4f39bed6999f1fc2987b39984b.de-de.xlf", "", "", "63a398a4-7e50-4038-a6ee-4d0ede53b8bc.3026d92881defd4f39bed6999f1fc2987b39984b.de-de.xlf") | Out-Null

$wsOverview.Select()
$wsOverview.Range("A1").Select()
